# Fix: the t-range for the first-order ODE (dx/dt = x, x(0) = 1) Euler
# approximation was changed from [0, 3] to [0, 2], keeping the same number
# of sample points (30 rows -> 29 steps). Recompute columns A (t) and B (x)
# for rows 3 through 31 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$t0 = 0.0
$t1 = 2.0
$n  = 29              # number of steps (30 points total, rows 2..31)
$dt = ($t1 - $t0) / $n

$t = $t0
$x = 1.0

for ($i = 0; $i -le $n; $i++) {
    $row = 2 + $i

    if ($i -gt 0) {
        $x = $x + $dt * $x
        $t = $t0 + $i * $dt
    }

    $ws.Cells.Item($row, 1).Value = $t
    $ws.Cells.Item($row, 2).Value = $x
}
